$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F: header "REX_DEF" matching the format of existing headers (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "REX_DEF"

# Add "[]" values for the new column's data rows (no special formatting, like C/E columns)
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
